$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 475
$ws.Range("I18").Value = 450
$ws.Range("K18").Value = 450
$ws.Range("M18").Value = -166
$ws.Range("H106").Value = 10446.167
$ws.Range("I106").Value = 6535.4
$ws.Range("K106").Value = 6535.4
$ws.Range("M106").Value = -5904.4
$ws.Range("H128").Value = 109999.5
$ws.Range("J128").Value = 109999.5
$ws.Range("L128").Value = 109999.5
$ws.Range("N128").Value = -119959.5
$ws.Range("H137").Value = 1491.0625
$ws.Range("I137").Value = 996.9231
$ws.Range("K137").Value = 2990.7693
$ws.Range("M137").Value = -440.7692999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24140.549
$ws.Range("I32").Value = 27169.715
$ws.Range("K32").Value = 27169.715
$ws.Range("M32").Value = -26882.715
$ws.Range("H61").Value = 1932.6364
$ws.Range("I61").Value = 1925.9
$ws.Range("K61").Value = 1925.9
$ws.Range("M61").Value = -1713.9
$ws.Range("H74").Value = 58887.555
$ws.Range("I74").Value = 92725.09
$ws.Range("J74").Value = 5714.2856
$ws.Range("K74").Value = 92725.09
$ws.Range("L74").Value = 5714.2856
$ws.Range("M74").Value = -91851.09
$ws.Range("N74").Value = -7462.2856
$ws.Range("H76").Value = 35088
$ws.Range("J76").Value = 35088
$ws.Range("L76").Value = 35088
$ws.Range("N76").Value = -35764
$ws.Range("H77").Value = 58887.555
$ws.Range("I77").Value = 92725.09
$ws.Range("J77").Value = 5714.2856
$ws.Range("K77").Value = 463625.45
$ws.Range("L77").Value = 28571.428
$ws.Range("M77").Value = -459257.45
$ws.Range("N77").Value = -37307.428
$ws.Range("H79").Value = 35088
$ws.Range("J79").Value = 35088
$ws.Range("L79").Value = 35088
$ws.Range("N79").Value = -37428
$ws.Range("H97").Value = 8702
$ws.Range("I97").Value = 11424
$ws.Range("J97").Value = 4813.4287
$ws.Range("K97").Value = 11424
$ws.Range("L97").Value = 4813.4287
$ws.Range("M97").Value = -10928
$ws.Range("N97").Value = -5805.4287
$ws.Range("H110").Value = 1589.2333
$ws.Range("I110").Value = 1671.6428
$ws.Range("K110").Value = 1671.6428
$ws.Range("M110").Value = 373.3571999999999
$ws.Range("H132").Value = 40917.77
$ws.Range("I132").Value = 45857.61
$ws.Range("K132").Value = 137572.83
$ws.Range("M132").Value = -135042.83
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 68000
$ws.Range("I134").Value = 68000
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 68000
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -62930
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 1932.6364
$ws.Range("I136").Value = 1925.9
$ws.Range("K136").Value = 5777.700000000001
$ws.Range("M136").Value = -3227.700000000001
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1341.75
$ws.Range("I19").Value = 1622.3334
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = 1622.3334
$ws.Range("L19").Value = 500
$ws.Range("M19").Value = -1452.3334
$ws.Range("N19").Value = -840
$ws.Range("H24").Value = 1341.75
$ws.Range("I24").Value = 1622.3334
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = 1622.3334
$ws.Range("L24").Value = 500
$ws.Range("M24").Value = -1452.3334
$ws.Range("N24").Value = -840
$ws.Range("H31").Value = 1793.1875
$ws.Range("I31").Value = 1783.4
$ws.Range("J31").Value = 1940
$ws.Range("K31").Value = 1783.4
$ws.Range("L31").Value = 1940
$ws.Range("M31").Value = -1488.4
$ws.Range("N31").Value = -2530
$ws.Range("H34").Value = 1793.1875
$ws.Range("I34").Value = 1783.4
$ws.Range("J34").Value = 1940
$ws.Range("K34").Value = 1783.4
$ws.Range("L34").Value = 1940
$ws.Range("M34").Value = -1581.4
$ws.Range("N34").Value = -2344
$ws.Range("H58").Value = 113477.22
$ws.Range("I58").Value = 144775.86
$ws.Range("J58").Value = 3932
$ws.Range("K58").Value = 144775.86
$ws.Range("L58").Value = 3932
$ws.Range("M58").Value = -144572.86
$ws.Range("N58").Value = -4338
$ws.Range("H136").Value = 113477.22
$ws.Range("I136").Value = 144775.86
$ws.Range("J136").Value = 3932
$ws.Range("K136").Value = 434327.58
$ws.Range("L136").Value = 11796
$ws.Range("M136").Value = -431777.58
$ws.Range("N136").Value = -16896

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 356.4074
$ws.Range("I92").Value = 352.375
$ws.Range("J92").Value = 362.27274
$ws.Range("K92").Value = 1057.125
$ws.Range("L92").Value = 1086.81822
$ws.Range("M92").Value = 190.875
$ws.Range("N92").Value = -3582.81822
$ws.Range("H137").Value = 3283.0476
$ws.Range("I137").Value = 1264.3334
$ws.Range("K137").Value = 3793.0002
$ws.Range("M137").Value = 1306.9998
$ws.Range("H141").Value = 1565.5
$ws.Range("I141").Value = 1565.5
$ws.Range("K141").Value = 4696.5
$ws.Range("M141").Value = 483.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 831.9
$ws.Range("I97").Value = 553.5
$ws.Range("J97").Value = 1249.5
$ws.Range("K97").Value = 553.5
$ws.Range("L97").Value = 1249.5
$ws.Range("M97").Value = -57.5
$ws.Range("N97").Value = -2241.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11904.529
$ws.Range("I7").Value = 27778.4
$ws.Range("J7").Value = 5290.4165
$ws.Range("K7").Value = 27778.4
$ws.Range("L7").Value = 5290.4165
$ws.Range("M7").Value = -27666.4
$ws.Range("N7").Value = -5514.4165
$ws.Range("H46").Value = 12690.469
$ws.Range("I46").Value = 17818.525
$ws.Range("K46").Value = 17818.525
$ws.Range("M46").Value = -17630.525
$ws.Range("H122").Value = 4154.963
$ws.Range("I122").Value = 3417.6365
$ws.Range("J122").Value = 4661.875
$ws.Range("K122").Value = 10252.9095
$ws.Range("L122").Value = 13985.625
$ws.Range("M122").Value = -7802.9095
$ws.Range("N122").Value = -18885.625
$ws.Range("H126").Value = 11904.529
$ws.Range("I126").Value = 27778.4
$ws.Range("J126").Value = 5290.4165
$ws.Range("K126").Value = 83335.20000000001
$ws.Range("L126").Value = 15871.2495
$ws.Range("M126").Value = -80865.20000000001
$ws.Range("N126").Value = -20811.2495

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1483.875
$ws.Range("I113").Value = 1515.4166
$ws.Range("J113").Value = 1452.3334
$ws.Range("K113").Value = 4546.2498
$ws.Range("L113").Value = 4357.0002
$ws.Range("M113").Value = -2376.2498
$ws.Range("N113").Value = -8697.0002
$ws.Range("H126").Value = 145314.84
$ws.Range("I126").Value = 207166.56
$ws.Range("K126").Value = 621499.6799999999
$ws.Range("M126").Value = -619029.68
$ws.Range("H132").Value = 79068.67999999999
$ws.Range("I132").Value = 85552.914
$ws.Range("K132").Value = 256658.742
$ws.Range("M132").Value = -254128.742
$ws.Range("H136").Value = 1874.4642
$ws.Range("I136").Value = 1857.3334
$ws.Range("K136").Value = 5572.0002
$ws.Range("M136").Value = -3022.0002
